# Rewrite speaker notes to concise reference style (9 slides in this deck).
# Each target slide's Notes Placeholder holds a single paragraph / single run;
# we just overwrite the TextRange.Text on that run's containing shape.

$p = $ppt.ActivePresentation

function Set-NotesText {
    param(
        [int]$SlideIndex,
        [string]$NewText
    )

    $slide = $p.Slides.Item($SlideIndex)
    $notes = $slide.NotesPage

    for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
        $shp = $notes.Shapes.Item($i)
        if ($shp.Name -eq "Notes Placeholder 2") {
            $shp.TextFrame.TextRange.Text = $NewText
            return
        }
    }

    throw "Notes Placeholder not found on slide $SlideIndex"
}

Set-NotesText 1 "Session 3. Requires Claude Code or equivalent agentic tool."

Set-NotesText 12 "The reveal. Agent found context beyond the explicit request."

Set-NotesText 13 "Transition to trust. AI has the keys — how to verify?"

Set-NotesText 14 "3-4 min. Group discussion: accuracy, trust vs verify. ‘Feels right’ = most dangerous moment."

Set-NotesText 15 "Core mindset: human reviews, judges, approves. Agent does grunt work."

Set-NotesText 2 "Strategist → Operator. AI talks and playbook done. Today: AI does the work."

Set-NotesText 3 "Callback to W2 cliffhanger. Audience still bottleneck: copy-paste, search, assemble."

Set-NotesText 4 "Copy-paste cycle. Every step has ‘you’ — human is bottleneck at each point."

Set-NotesText 6 "Same four steps. Human only at step 4 (review). Agent handles retrieve, identify, create."
